$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPrompt = "Please rewrite this in the style of Jason Fladlien:"

$values = @{
    2  = "185.74622321128845"
    3  = "87.75262379646301"
    4  = "61.15979361534119"
    5  = "141.45266103744507"
    6  = "154.113951921463"
    7  = "143.74763584136963"
    8  = "77.22561883926392"
    9  = "68.87961220741272"
    10 = "161.40382480621338"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 10)
    $cell.Value = "$newPrompt`n$($values[$row])"
}
